$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H75").Value = 41638.25
$ws.Range("J75").Value = 41638.25
$ws.Range("L75").Value = 41638.25
$ws.Range("N75").Value = -43510.25

$ws.Range("H78").Value = 41638.25
$ws.Range("J78").Value = 41638.25
$ws.Range("L78").Value = 124914.75
$ws.Range("N78").Value = -134274.75

$ws.Range("H141").Value = 5075.9375
$ws.Range("I141").Value = 3240.9092
$ws.Range("K141").Value = 9722.7276
$ws.Range("M141").Value = -4542.7276

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 35673.38
$ws.Range("I32").Value = 35562.066
$ws.Range("J32").Value = 36953.5
$ws.Range("K32").Value = 35562.066
$ws.Range("L32").Value = 36953.5
$ws.Range("M32").Value = -35275.066
$ws.Range("N32").Value = -37527.5

$ws.Range("H80").Value = 50000
$ws.Range("J80").Value = 50000
$ws.Range("L80").Value = 50000
$ws.Range("N80").Value = -51996

$ws.Range("H83").Value = 50000
$ws.Range("J83").Value = 50000
$ws.Range("L83").Value = 150000
$ws.Range("N83").Value = -159984

$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 46401
$ws.Range("J92").Value = 46401
$ws.Range("L92").Value = 46401
$ws.Range("N92").Value = -51393

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H124").Value = 29794.666
$ws.Range("J124").Value = 29794.666
$ws.Range("L124").Value = 29794.666
$ws.Range("N124").Value = -34704.666

$ws.Range("H131").Value = 38318
$ws.Range("J131").Value = 38318
$ws.Range("L131").Value = 38318
$ws.Range("N131").Value = -48398

$ws.Range("H141").Value = 13396.4
$ws.Range("J141").Value = 14998
$ws.Range("L141").Value = 14998
$ws.Range("N141").Value = -25358

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3082.1667
$ws.Range("I5").Value = 8568.833000000001
$ws.Range("J5").Value = 1253.2778
$ws.Range("K5").Value = 25706.499
$ws.Range("L5").Value = 3759.8334
$ws.Range("M5").Value = -25594.499
$ws.Range("N5").Value = -3983.8334

$ws.Range("H44").Value = 189006.75
$ws.Range("I44").Value = 500816.66
$ws.Range("J44").Value = 1920.8
$ws.Range("K44").Value = 1502449.98
$ws.Range("L44").Value = 5762.4
$ws.Range("M44").Value = -1502051.98
$ws.Range("N44").Value = -6558.4

$ws.Range("H48").Value = 200201340
$ws.Range("J48").Value = 200201340
$ws.Range("L48").Value = 600604020
$ws.Range("N48").Value = -600604520

$ws.Range("H98").Value = 913.0909
$ws.Range("I98").Value = 744.3333
$ws.Range("J98").Value = 976.375
$ws.Range("K98").Value = 2232.9999
$ws.Range("L98").Value = 2929.125
$ws.Range("M98").Value = -734.9998999999998
$ws.Range("N98").Value = -5925.125

$ws.Range("H113").Value = 4509.4614
$ws.Range("I113").Value = 5871.1055
$ws.Range("J113").Value = 813.5714
$ws.Range("K113").Value = 17613.3165
$ws.Range("L113").Value = 2440.7142
$ws.Range("M113").Value = -15443.3165
$ws.Range("N113").Value = -6780.7142

$ws.Range("H131").Value = 48911.09
$ws.Range("J131").Value = 56863.418
$ws.Range("L131").Value = 170590.254
$ws.Range("N131").Value = -180670.254

$ws.Range("H135").Value = 3082.1667
$ws.Range("I135").Value = 8568.833000000001
$ws.Range("J135").Value = 1253.2778
$ws.Range("K135").Value = 77119.497
$ws.Range("L135").Value = 11279.5002
$ws.Range("M135").Value = -74584.497
$ws.Range("N135").Value = -16349.5002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 31497.334
$ws.Range("J74").Value = 33796.8
$ws.Range("L74").Value = 33796.8
$ws.Range("N74").Value = -35668.8

$ws.Range("H77").Value = 31497.334
$ws.Range("J77").Value = 33796.8
$ws.Range("L77").Value = 101390.4
$ws.Range("N77").Value = -110750.4

$ws.Range("H118").Value = 28870.666
$ws.Range("J118").Value = 28870.666
$ws.Range("L118").Value = 28870.666
$ws.Range("N118").Value = -32184.666

$ws.Range("H127").Value = 43438.668
$ws.Range("J127").Value = 43438.668
$ws.Range("L127").Value = 43438.668
$ws.Range("N127").Value = -53358.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H88").Value = 44189
$ws.Range("J88").Value = 44189
$ws.Range("L88").Value = 44189
$ws.Range("N88").Value = -45045

$ws.Range("H91").Value = 44189
$ws.Range("J91").Value = 44189
$ws.Range("L91").Value = 44189
$ws.Range("N91").Value = -47153

$ws.Range("H102").Value = 48553
$ws.Range("J102").Value = 48553
$ws.Range("L102").Value = 48553
$ws.Range("N102").Value = -55043

$ws.Range("H109").Value = 35281
$ws.Range("J109").Value = 35281
$ws.Range("L109").Value = 35281
$ws.Range("N109").Value = -38055

$ws.Range("H123").Value = 32872.668
$ws.Range("J123").Value = 32872.668
$ws.Range("L123").Value = 32872.668
$ws.Range("N123").Value = -42672.668

$ws.Range("H129").Value = 31904.5
$ws.Range("J129").Value = 31904.5
$ws.Range("L129").Value = 31904.5
$ws.Range("N129").Value = -41904.5

$ws.Range("H131").Value = 43326
$ws.Range("J131").Value = 43326
$ws.Range("L131").Value = 43326
$ws.Range("N131").Value = -53406

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 26525.2
$ws.Range("J27").Value = 26525.2
$ws.Range("L27").Value = 26525.2
$ws.Range("N27").Value = -26663.2

$ws.Range("H69").Value = 23923.5
$ws.Range("J69").Value = 23923.5
$ws.Range("L69").Value = 23923.5
$ws.Range("N69").Value = -25421.5

$ws.Range("H72").Value = 23923.5
$ws.Range("J72").Value = 23923.5
$ws.Range("L72").Value = 71770.5
$ws.Range("N72").Value = -79258.5

$ws.Range("H80").Value = 35825
$ws.Range("J80").Value = 35825
$ws.Range("L80").Value = 35825
$ws.Range("N80").Value = -37821

$ws.Range("H83").Value = 35825
$ws.Range("J83").Value = 35825
$ws.Range("L83").Value = 107475
$ws.Range("N83").Value = -117459

$ws.Range("H92").Value = 64440
$ws.Range("J92").Value = 64440
$ws.Range("L92").Value = 64440
$ws.Range("N92").Value = -69432

$ws.Range("H103").Value = 48574.332
$ws.Range("J103").Value = 48574.332
$ws.Range("L103").Value = 48574.332
$ws.Range("N103").Value = -50918.332

$ws.Range("H109").Value = 36251.332
$ws.Range("J109").Value = 36251.332
$ws.Range("L109").Value = 36251.332
$ws.Range("N109").Value = -39025.332

$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

$ws.Range("H118").Value = 29696
$ws.Range("J118").Value = 44392
$ws.Range("L118").Value = 44392
$ws.Range("N118").Value = -47706

$ws.Range("H127").Value = 31561.334
$ws.Range("J127").Value = 31561.334
$ws.Range("L127").Value = 31561.334
$ws.Range("N127").Value = -41481.334
